$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set up the new column (AS = column 45) width to match the other
#     data columns (OOXML width 12, which corresponds to ColumnWidth 11.17
#     in Excel's character-count units for this font/theme). ---
$ws.Columns.Item(45).ColumnWidth = 11.17

# --- Style template cells: reuse existing formatting from column AR (44)
#     so the new column's cells land on the exact same style records
#     (s="1" normal, s="2" yellow <125, s="3" light-blue 125-139.9)
#     instead of Excel fabricating brand-new style entries. ---
$styleSrc1 = $ws.Cells.Item(2, 44)   # style 1 : no fill  (value >= 140)
$styleSrc2 = $ws.Cells.Item(18, 44)  # style 2 : yellow   (value < 125)
$styleSrc3 = $ws.Cells.Item(9, 44)   # style 3 : light blue (125 <= value < 140)

# --- Header cell AS1: "2024/10/23" ---
# A bare $cell.Value = "2024/10/23" gets auto-parsed by Excel as a date
# serial (like real Excel does for General-formatted cells). Enter it
# quote-prefixed so it is kept as literal text, then repaint the cell's
# format (font/fill/number format) from AR1 so it ends up on the same
# style record as the rest of the header row.
$styleSrc1.Copy($ws.Cells.Item(1, 45))
$ws.Cells.Item(1, 45).Value = "'2024/10/23"
$ws.Cells.Item(1, 44).Copy()
$ws.Cells.Item(1, 45).PasteSpecial(-4122)

$styleSrc3.Copy($ws.Cells.Item(2, 45))
$ws.Cells.Item(2, 45).Value = 137.1
$styleSrc1.Copy($ws.Cells.Item(3, 45))
$ws.Cells.Item(3, 45).Value = 165.1
$styleSrc1.Copy($ws.Cells.Item(4, 45))
$ws.Cells.Item(4, 45).Value = 213.8
$styleSrc1.Copy($ws.Cells.Item(5, 45))
$ws.Cells.Item(5, 45).Value = 177.8
$styleSrc1.Copy($ws.Cells.Item(6, 45))
$ws.Cells.Item(6, 45).Value = 160.7
$styleSrc1.Copy($ws.Cells.Item(7, 45))
$ws.Cells.Item(7, 45).Value = 189.2
$styleSrc1.Copy($ws.Cells.Item(8, 45))
$ws.Cells.Item(8, 45).Value = 152.1
$styleSrc2.Copy($ws.Cells.Item(9, 45))
$ws.Cells.Item(9, 45).Value = 116.8
$styleSrc1.Copy($ws.Cells.Item(10, 45))
$ws.Cells.Item(10, 45).Value = 153.4
$styleSrc1.Copy($ws.Cells.Item(11, 45))
$ws.Cells.Item(11, 45).Value = 154.8
$styleSrc1.Copy($ws.Cells.Item(12, 45))
$ws.Cells.Item(12, 45).Value = 168
$styleSrc2.Copy($ws.Cells.Item(13, 45))
$ws.Cells.Item(13, 45).Value = 115.6
$styleSrc1.Copy($ws.Cells.Item(14, 45))
$ws.Cells.Item(14, 45).Value = 290.6
$styleSrc2.Copy($ws.Cells.Item(15, 45))
$ws.Cells.Item(15, 45).Value = 109.4
$styleSrc1.Copy($ws.Cells.Item(16, 45))
$ws.Cells.Item(16, 45).Value = 160
$styleSrc3.Copy($ws.Cells.Item(17, 45))
$ws.Cells.Item(17, 45).Value = 129.1
$styleSrc1.Copy($ws.Cells.Item(18, 45))
$ws.Cells.Item(18, 45).Value = 164.4
$styleSrc1.Copy($ws.Cells.Item(19, 45))
$ws.Cells.Item(19, 45).Value = 144.3
$styleSrc1.Copy($ws.Cells.Item(20, 45))
$ws.Cells.Item(20, 45).Value = 154.3
$styleSrc1.Copy($ws.Cells.Item(21, 45))
$ws.Cells.Item(21, 45).Value = 177.9
$styleSrc1.Copy($ws.Cells.Item(22, 45))
$ws.Cells.Item(22, 45).Value = 178.9
$styleSrc1.Copy($ws.Cells.Item(23, 45))
$ws.Cells.Item(23, 45).Value = 193.1
$styleSrc1.Copy($ws.Cells.Item(24, 45))
$ws.Cells.Item(24, 45).Value = 144.5
$styleSrc1.Copy($ws.Cells.Item(25, 45))
$ws.Cells.Item(25, 45).Value = 165.3
$styleSrc1.Copy($ws.Cells.Item(26, 45))
$ws.Cells.Item(26, 45).Value = 177.5
$styleSrc1.Copy($ws.Cells.Item(27, 45))
$ws.Cells.Item(27, 45).Value = 149.1
$styleSrc3.Copy($ws.Cells.Item(28, 45))
$ws.Cells.Item(28, 45).Value = 133.8
$styleSrc1.Copy($ws.Cells.Item(29, 45))
$ws.Cells.Item(29, 45).Value = 170.7
$styleSrc3.Copy($ws.Cells.Item(30, 45))
$ws.Cells.Item(30, 45).Value = 128.4
$styleSrc1.Copy($ws.Cells.Item(31, 45))
$ws.Cells.Item(31, 45).Value = 344.4
$styleSrc1.Copy($ws.Cells.Item(32, 45))
$ws.Cells.Item(32, 45).Value = 156.8
$styleSrc1.Copy($ws.Cells.Item(33, 45))
$ws.Cells.Item(33, 45).Value = 283.4
$styleSrc2.Copy($ws.Cells.Item(34, 45))
$ws.Cells.Item(34, 45).Value = 120.8
$styleSrc1.Copy($ws.Cells.Item(35, 45))
$ws.Cells.Item(35, 45).Value = 184.2
$styleSrc3.Copy($ws.Cells.Item(36, 45))
$ws.Cells.Item(36, 45).Value = 134.6
$styleSrc1.Copy($ws.Cells.Item(37, 45))
$ws.Cells.Item(37, 45).Value = 143.9
$styleSrc1.Copy($ws.Cells.Item(38, 45))
$ws.Cells.Item(38, 45).Value = 156.4
$styleSrc1.Copy($ws.Cells.Item(39, 45))
$ws.Cells.Item(39, 45).Value = 141.3
$styleSrc1.Copy($ws.Cells.Item(40, 45))
$ws.Cells.Item(40, 45).Value = 206.4
$styleSrc1.Copy($ws.Cells.Item(41, 45))
$ws.Cells.Item(41, 45).Value = 150.3
$styleSrc1.Copy($ws.Cells.Item(42, 45))
$ws.Cells.Item(42, 45).Value = 180.6
$styleSrc1.Copy($ws.Cells.Item(43, 45))
$ws.Cells.Item(43, 45).Value = 158.8
$styleSrc1.Copy($ws.Cells.Item(44, 45))
$ws.Cells.Item(44, 45).Value = 166.2
$styleSrc3.Copy($ws.Cells.Item(45, 45))
$ws.Cells.Item(45, 45).Value = 130.2
$styleSrc1.Copy($ws.Cells.Item(46, 45))
$ws.Cells.Item(46, 45).Value = 151.1
$styleSrc3.Copy($ws.Cells.Item(47, 45))
$ws.Cells.Item(47, 45).Value = 133.7
$styleSrc1.Copy($ws.Cells.Item(48, 45))
$ws.Cells.Item(48, 45).Value = 167
$styleSrc3.Copy($ws.Cells.Item(49, 45))
$ws.Cells.Item(49, 45).Value = 135
$styleSrc1.Copy($ws.Cells.Item(50, 45))
$ws.Cells.Item(50, 45).Value = 144.1
$styleSrc1.Copy($ws.Cells.Item(51, 45))
$ws.Cells.Item(51, 45).Value = 166.2
$styleSrc3.Copy($ws.Cells.Item(52, 45))
$ws.Cells.Item(52, 45).Value = 138.1
$styleSrc1.Copy($ws.Cells.Item(53, 45))
$ws.Cells.Item(53, 45).Value = 149.4

$excel.CutCopyMode = $false

